$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column "priority" (I) cleanup -----------------------------------
# Only the America/New_York rows (CN=12, US=59) keep a priority value;
# everything else loses its stray priority number.
$ws.Range("I4").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("I15").ClearContents()
$ws.Range("I17:I48").ClearContents()
$ws.Range("I51").ClearContents()
$ws.Range("I52").ClearContents()
$ws.Range("I53").ClearContents()
$ws.Range("I54").ClearContents()
$ws.Range("I57").ClearContents()
$ws.Range("I61").ClearContents()
$ws.Range("I64").ClearContents()

$ws.Range("I12").Value = 2
$ws.Range("I59").Value = 1

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 24.7109375
$ws.Columns.Item(7).ColumnWidth = 10.28515625

# --- AutoFilter on the timezone column (D), showing only America/* ----
$vals = @("America/Argentina/Buenos_Aires","America/Caracas","America/Cuiaba","America/Mexico_City","America/New_York","America/Santiago","America/Sao_Paulo","America/Toronto")
$ws.Range("A1:H64").AutoFilter(4, $vals, 7)

# --- Selection moved to A59 --------------------------------------------
$ws.Range("A59").Select()

# --- Shrink the workbook window width (matches the recorded view state)
$excel.ActiveWindow.Width = 17265
